# Commit: "add bottom temp indicator"
#
# Row 6 of the indicator table was a placeholder row for the "Bottom
# temperature" indicator. Fill it in with the real 2024 status, a more
# precise indicator label, and the actual time-series image filename,
# and bump the image height factor (column F, "h") to 1.5 to fit the
# new chart.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A6: Indicator (units) - clarify region covered by the bottom temp data
$ws.Range("A6").Value = "Bottom temperature in MAB and SNE(°C) "

# B6: Status in 2024 - replace placeholder text with the real status
$ws.Range("B6").Value = "Above long term average (Fall); near long term average (Spring)"

# D6: Time series image filename - point at the real chart export
$ws.Range("D6").Value = "BottomT_2025-04-10.png"

# F6: image height factor bumped from 1 to 1.5
$ws.Range("F6").Value = 1.5

# Leave the selection on C6, matching where the author's cursor ended up
$ws.Range("C6").Select()
